$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to Text so numeric-looking strings
# (e.g. "1.001") are stored verbatim instead of being parsed as numbers,
# matching the existing inline-string cells in this sheet.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.479.14"
$ws.Range("E2").Value = "  -0.84%  "

$ws.Range("D3").Value = "1.893.75"
$ws.Range("E3").Value = "  +0.16%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.19%  "

$ws.Range("D5").Value = "237.96"
$ws.Range("E5").Value = "  +0.73%  "

$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.10%  "

$ws.Range("D7").Value = "0.4897"
$ws.Range("E7").Value = "  +0.29%  "

$ws.Range("D8").Value = "0.2932"
$ws.Range("E8").Value = "  +1.08%  "

$ws.Range("D9").Value = "0.06685"
$ws.Range("E9").Value = "  +0.22%  "

$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").Value = "17.08"
$ws.Range("E10").Value = "  +2.37%  "

$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").Value = "1.848.76"
$ws.Range("E11").Value = "  -1.96%  "

$ws.Range("D12").Value = "0.07349"
$ws.Range("E12").Value = "  +1.39%  "

$ws.Range("D13").Value = "5.145"
$ws.Range("E13").Value = "  +2.68%  "

$ws.Range("D14").Value = "88.09"
$ws.Range("E14").Value = "  -1.38%  "

$ws.Range("D15").Value = "0.6649"
$ws.Range("E15").Value = "  +0.02%  "

$ws.Range("D16").Value = "30.437.56"
$ws.Range("E16").Value = "  -0.76%  "

$ws.Range("D17").Value = "13.46"
$ws.Range("E17").Value = "  +3.68%  "

$ws.Range("D18").Value = "0.000007830"
$ws.Range("E18").Value = "  -0.52%  "

$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  +0.00%  "

$ws.Range("D20").Value = "2.129.25"
$ws.Range("E20").Value = "  +0.05%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "5.293"
$ws.Range("E21").Value = "  +11.72%  "

$ws.Range("B22").Value = "BinanceUSD"
$ws.Range("C22").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.34%  "

$ws.Range("D23").Value = "188.90"
$ws.Range("E23").Value = "  -0.79%  "

$ws.Range("D24").Value = "6.157"

$ws.Range("D25").Value = "9.483"
$ws.Range("E25").Value = "  +2.02%  "

$ws.Range("D26").Value = "164.70"
$ws.Range("E26").Value = "  +3.25%  "

$ws.Range("D27").Value = "18.32"
$ws.Range("E27").Value = "  -0.11%  "

$ws.Range("D28").Value = "1.931"
$ws.Range("E28").Value = "  +5.76%  "

$ws.Range("D29").Value = "1.465"
$ws.Range("E29").Value = "  +4.27%  "

$ws.Range("D30").Value = "4.352"
$ws.Range("E30").Value = "  +2.15%  "

$ws.Range("D31").Value = "0.09181"
$ws.Range("E31").Value = "  +1.61%  "

$ws.Range("D32").Value = "4.082"
$ws.Range("E32").Value = "  +3.54%  "

$ws.Range("D33").Value = "0.05207"
$ws.Range("E33").Value = "  +0.05%  "

$ws.Range("D34").Value = "0.7412"
$ws.Range("E34").Value = "  +1.29%  "

$ws.Range("D35").Value = "1.099"
$ws.Range("E35").Value = "  +1.35%  "

$ws.Range("D36").Value = "2.717"
$ws.Range("E36").Value = "  +0.84%  "

$ws.Range("D37").Value = "0.01817"
$ws.Range("E37").Value = "  -0.43%  "

$ws.Range("D38").Value = "2.672"
$ws.Range("E38").Value = "  +0.02%  "

$ws.Range("D39").Value = "0.9163"
$ws.Range("E39").Value = "  -0.71%  "

$ws.Range("D40").Value = "2.041"
$ws.Range("E40").Value = "  -0.88%  "

$ws.Range("D41").Value = "0.4400"
$ws.Range("E41").Value = "  -1.38%  "

$ws.Range("D42").Value = "5.944"
$ws.Range("E42").Value = "  +3.91%  "

$ws.Range("D43").Value = "105.91"
$ws.Range("E43").Value = "  +1.29%  "

$ws.Range("D44").Value = "0.9926"
$ws.Range("E44").Value = "  -0.70%  "

$ws.Range("D45").Value = "0.1387"
$ws.Range("E45").Value = "  +3.82%  "

$ws.Range("E46").Value = "  +18.62%  "

$ws.Range("D47").Value = "7.569"
$ws.Range("E47").Value = "  +3.60%  "

$ws.Range("D48").Value = "9.022"
$ws.Range("E48").Value = "  +4.72%  "

$ws.Range("D49").Value = "34.93"
$ws.Range("E49").Value = "  +4.98%  "

$ws.Range("D50").Value = "0.05826"
$ws.Range("E50").Value = "  -0.13%  "

$ws.Range("D51").Value = "0.3947"
$ws.Range("E51").Value = "  -8.13%  "

# Restore default (General) formatting so no stray style indices are
# introduced on cells beyond the inline-string value itself.
$ws.Range("D2:E51").ClearFormats()

